$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column (H), matching the look of the existing
# header row (bold, thin box border, centered/top-aligned - same as
# B1:G1).
$ws.Range("H1").Value = "Save"
$ws.Range("H1").Font.Bold = $true
$ws.Range("H1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("H1").VerticalAlignment = -4160     # xlTop
$ws.Range("H1").Borders.LineStyle = 1         # xlContinuous (thin box)

# New data values for the Save column
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
